$d = $word.ActiveDocument

# --- Change 1: insert a new empty (bold, Calibri 12pt, en-GB) paragraph
# right before the "EXTRA FEATURES NOT PERTAINING..." paragraph, duplicating
# the formatting of the existing empty paragraph that precedes it.
$marker = "EXTRA FEATURES NOT PERTAINING AN .xml OR .java FILE/NOTES ON EXTRA FEATURES ABOVE (see checklist for feature number):"
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$marker*") {
        $p.Range.InsertParagraphBefore()
        break
    }
}

# --- Change 2: update the wording of the "about screen" bullet (#8) to
# describe the LinkedIn / app-opening fix.
$old = ' Added an ' + [char]8220 + 'about' + [char]8221 + ' screen with four clickable TextViews. The contact one opens the default emailing app and composes an email directed to my uni address. All the other ones below redirect to the app' + [char]8217 + 's repository, my Github profile, and my LinkedIn profile, either in the default browser or in their respective apps. Unfortunately this feature doesn' + [char]8217 + 't seem to work on emulated devices, but it works on my Huawei Nova 5 Pro flawlessly, and is the device I used for my demo.'
$new = ' Added an ' + [char]8220 + 'about' + [char]8221 + ' screen with four functional TextViews (the rest of TextViews are also clickable and focusable, but don' + [char]8217 + 't redirect to anything). The contact one opens the default emailing app and composes an email directed to my uni address. All the other ones below redirect to the app' + [char]8217 + 's repository, my Github profile, and my LinkedIn profile, either in the default browser or in their respective apps. Unfortunately this feature doesn' + [char]8217 + 't seem to work on emulated devices, but it works on my Huawei Nova 5 Pro flawlessly, and is the device I used for my demo.'

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
